# Fix test case "Log What's Your Name" (cl-ui TestGrid.xlsx)
# - Splits the "Close/Delete -> Save Branch/Abandon Branch" button-label
#   assertion out into its own explicit test row.
# - Splits "After branch is saved the original branch should retain all of
#   its original training" into two separate manual-only checks: one for
#   saving a branch, one for abandoning a branch.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Insert two new rows inside the "Branching" block of the table -------
# Row 26: new row asserting the buttons change label after branching.
$ws.Rows.Item(26).Insert()
# Row 32 (post first insert): new row for the "abandon branch" unchanged-training check.
$ws.Rows.Item(32).Insert()

# Grow the table/autofilter/dimension to cover the two newly inserted rows.
$tbl.Resize($ws.Range("A1:F59"))

# --- Populate the first new row (26) --------------------------------------
$ws.Range("A26").Value = "Train"
$ws.Range("B26").Value = "Branching"
$ws.Range("C26").Value = "Train Dialog"
$ws.Range("D26").Value = '"Close" and "Delete" buttons should change to "Save Branch" and "Abandon Branch" after branching'

# --- Rework the old "retain all of its original training" row (now 31) ---
# into the "saved" variant, and mark it manual.
$ws.Range("D31").Value = "After branch is saved the original training should remain unchanged"
$ws.Range("E31").Value = "manual"

# --- Populate the second new row (32) with the "abandoned" variant -------
$ws.Range("A32").Value = "Train"
$ws.Range("B32").Value = "Branching"
$ws.Range("C32").Value = "Train Dialog"
$ws.Range("D32").Value = "After branch is abandonded the original training should remain unchanged"
$ws.Range("E32").Value = "manual"

# --- Update the current selection/active cell to reflect the new layout ---
$ws.Range("D33").Select()
